$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "last updated" timestamp (15:38 -> 16:55) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 16:55"

# --- COVID-19 stats refresh for a batch of countries (B=Casos totales, C=Nuevos casos,
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6260879
$ws.Range("C4").Value = 3308
$ws.Range("D4").Value = 3498208
$ws.Range("E4").Value = 2573696
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = 188975

# Row 6: India
$ws.Range("B6").Value = 3810625
$ws.Range("C6").Value = 44517
$ws.Range("D6").Value = 2931005
$ws.Range("E6").Value = 812749
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 411
$ws.Range("H6").Value = 66871

# Row 18: Arabia Saudita
$ws.Range("B18").Value = 317486
$ws.Range("C18").Value = 816
$ws.Range("D18").Value = 292510
$ws.Range("E18").Value = 21020
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = 3956

# Row 23: Alemania
$ws.Range("B23").Value = 246499
$ws.Range("C23").Value = 498
$ws.Range("D23").Value = 221800
$ws.Range("E23").Value = 15314
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 9385

# Row 24: Irak
$ws.Range("B24").Value = 242284
$ws.Range("C24").Value = 3946
$ws.Range("D24").Value = 184205
$ws.Range("E24").Value = 50878
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 78
$ws.Range("H24").Value = 7201

# Row 62: Uzbekistan
$ws.Range("B62").Value = 42437
$ws.Range("C62").Value = 310
$ws.Range("D62").Value = 39664
$ws.Range("E62").Value = 2446
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 327

# Row 66: Moldavia
$ws.Range("B66").Value = 37740
$ws.Range("C66").Value = 532
$ws.Range("D66").Value = 26189
$ws.Range("E66").Value = 10527
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 16
$ws.Range("H66").Value = 1024

# Row 68: Kenia
$ws.Range("B68").Value = 34493
$ws.Range("C68").Value = 178
$ws.Range("D68").Value = 20211
$ws.Range("E68").Value = 13701
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 581

# Row 69: Serbia
$ws.Range("B69").Value = 31581
$ws.Range("C69").Value = 99
$ws.Range("D69").Value = 30188
$ws.Range("E69").Value = 677
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 716

# Row 89: Zambia
$ws.Range("B89").Value = 12415
$ws.Range("C89").Value = 34
$ws.Range("D89").Value = 11494
$ws.Range("E89").Value = 629
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 292

# Row 90: Noruega
$ws.Range("B90").Value = 10951
$ws.Range("C90").Value = 80
$ws.Range("D90").Value = 9348
$ws.Range("E90").Value = 1339
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 264

# Row 98: Tayikistan
$ws.Range("B98").Value = 8654
$ws.Range("C98").Value = 35
$ws.Range("D98").Value = 7447
$ws.Range("E98").Value = 1138
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 69

# Row 100: Haiti
$ws.Range("B100").Value = 8258
$ws.Range("C100").Value = 28
$ws.Range("D100").Value = 5870
$ws.Range("E100").Value = 2182
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 206

# --- Reunion overtakes Nueva Zelanda and Botsuana in the ranking (rows 148-150) ---
# Row 148 becomes Reunion with its updated figures
$ws.Range("A148").Value = "Reunion"
$ws.Range("B148").Value = 1796
$ws.Range("C148").Value = 82
$ws.Range("D148").Value = 880
$ws.Range("E148").Value = 906
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 10

# Row 149 becomes Nueva Zelanda (figures unchanged, just shifted down a row)
$ws.Range("A149").Value = "Nueva Zelanda"
$ws.Range("B149").Value = 1757
$ws.Range("C149").Value = 5
$ws.Range("D149").Value = 1606
$ws.Range("E149").Value = 129
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 22

# Row 150 becomes Botsuana (figures unchanged, just shifted down a row)
$ws.Range("A150").Value = "Botsuana"
$ws.Range("B150").Value = 1724
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 493
$ws.Range("E150").Value = 1225
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 6
